# Apply the edits described in the commit "test elec sec US files" to the
# Capacity Supply Curve workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "About" sheet: drop the state label ("Nevada") and the "last
#    updated" date that used to sit next to the title in row 1.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
[void]$wsAbout.Range("B1:C1").Clear()

# ---------------------------------------------------------------------
# 2) "CSC-CSCCCMvSoECBtY" sheet: just a saved-selection change.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("CSC-CSCCCMvSoECBtY")
[void]$ws2.Activate()
[void]$ws2.Range("B1:N2").Select()

# ---------------------------------------------------------------------
# 3) "CSC-CSCSoCECBiaSY" sheet: the share-of-cost-effective-capacity
#    values for year rows 7..30 (sheet rows 2..25) change from 1 to
#    0.3, except for row index 12 (sheet row 7) which becomes 0.2.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
[void]$ws4.Activate()

$ws4.Range("B2:AE6").Value = 0.3
$ws4.Range("B7:AE7").Value = 0.2
$ws4.Range("B8:AE25").Value = 0.3

[void]$ws4.Range("B7:AE7").Select()

# ---------------------------------------------------------------------
# 4) Re-activate the "About" sheet so it becomes the tab that is shown
#    / selected when the workbook is reopened.
# ---------------------------------------------------------------------
[void]$wsAbout.Activate()
